$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.013.65"
$ws.Range("E2").Value = "  -2.22%  "

# Row 3
$ws.Range("D3").Value = "2.574.66"
$ws.Range("E3").Value = "  -2.20%  "

# Row 4
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.28"
$ws.Range("E5").Value = "  +0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.14"
$ws.Range("E6").Value = "  -3.12%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  -1.48%  "

# Row 9
$ws.Range("D9").Value = "2.579.53"
$ws.Range("E9").Value = "  -2.79%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.63"
$ws.Range("E10").Value = "  -3.56%  "

# Row 11
$ws.Range("E11").Value = "  +0.11%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.162"
$ws.Range("E12").Value = "  +12.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.345"
$ws.Range("E13").Value = "  +1.42%  "

# Row 14
$ws.Range("D14").Value = "3.032.27"
$ws.Range("E14").Value = "  -2.02%  "

# Row 15
$ws.Range("D15").Value = "59.015.05"
$ws.Range("E15").Value = "  -2.02%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.41"
$ws.Range("E16").Value = "  +4.60%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").Value = "  +0.70%  "

# Row 18
$ws.Range("D18").Value = "2.575.22"
$ws.Range("E18").Value = "  -2.65%  "

# Row 19
$ws.Range("E19").Value = "  -0.23%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.85"
$ws.Range("E20").Value = "  -1.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.29"
$ws.Range("E21").Value = "  -0.32%  "

# Row 22
$ws.Range("E22").Value = "  -0.07%  "

# Row 23
$ws.Range("E23").Value = "  +0.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.27"
$ws.Range("E24").Value = "  -3.97%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.463"
$ws.Range("E25").Value = "  +6.88%  "

# Row 26
$ws.Range("E26").Value = "  +0.55%  "

# Row 27
$ws.Range("E27").Value = "  -2.58%  "

# Row 28
$ws.Range("E28").Value = "  -0.59%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0772"
$ws.Range("E29").Value = "  -0.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.06%  "

# Row 31
$ws.Range("E31").Value = "  -1.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.94"
$ws.Range("E32").Value = "  +2.59%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.05"
$ws.Range("E33").Value = "  -1.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.90"
$ws.Range("E34").Value = "  -1.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.99"
$ws.Range("E35").Value = "  -1.94%  "

# Row 36
$ws.Range("E36").Value = "  -1.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.869"
$ws.Range("E37").Value = "  -4.66%  "

# Row 38
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.47"
$ws.Range("E38").Value = "  +0.20%  "

# Row 39
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.870"
$ws.Range("E39").Value = "  -4.47%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.48"
$ws.Range("E40").Value = "  -1.97%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "295.59"
$ws.Range("E41").Value = "  -0.95%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.65"
$ws.Range("E42").Value = "  -0.36%  "

# Row 43
$ws.Range("E43").Value = "  +0.21%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.54"
$ws.Range("E44").Value = "  +9.82%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0973"
$ws.Range("E45").Value = "  +0.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.591"
$ws.Range("E46").Value = "  -2.43%  "

# Row 47
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.64"
$ws.Range("E47").Value = "  +0.12%  "

# Row 48
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0534"
$ws.Range("E48").Value = "  -2.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.11"
$ws.Range("E49").Value = "  -1.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0231"
$ws.Range("E50").Value = "  -0.29%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.44"
$ws.Range("E51").Value = "  -0.18%  "

